# Natmi following Dr Hou advice:
# "FAPs" is added as an additional sending cluster (alongside "sCs") for the
# Fgf15->Fgfr1 edge, and the specificity/weight metrics for every
# sCs-sourced row are recomputed now that FAPs participates too. This
# doubles the row count (5 -> 10) and grows the used range to A1:T11.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Fgf15"
$ws.Range("C2").Value = "Fgfr1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.029478
$ws.Range("H2").Value = 0.088434
$ws.Range("I2").Value = 0.1535387136874709
$ws.Range("J2").Value = 0.1535387136874709
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 10.25548833333333
$ws.Range("N2").Value = 30.766465
$ws.Range("O2").Value = 0.09996462463766297
$ws.Range("P2").Value = 0.09996462463766294
$ws.Range("Q2").Value = 0.30231128509
$ws.Range("R2").Value = 2.72080156581
$ws.Range("S2").Value = 0.01534843988111763
$ws.Range("T2").Value = 0.01534843988111763

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Fgf15"
$ws.Range("C3").Value = "Fgfr1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.029478
$ws.Range("H3").Value = 0.088434
$ws.Range("I3").Value = 0.1535387136874709
$ws.Range("J3").Value = 0.1535387136874709
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 74.64939600000001
$ws.Range("N3").Value = 223.948188
$ws.Range("O3").Value = 0.7276395436298834
$ws.Range("P3").Value = 0.7276395436298831
$ws.Range("Q3").Value = 2.200514895288
$ws.Range("R3").Value = 19.804634057592
$ws.Range("S3").Value = 0.1117208395570707
$ws.Range("T3").Value = 0.1117208395570706

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Fgf15"
$ws.Range("C4").Value = "Fgfr1"
$ws.Range("D4").Value = "M1"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.029478
$ws.Range("H4").Value = 0.088434
$ws.Range("I4").Value = 0.1535387136874709
$ws.Range("J4").Value = 0.1535387136874709
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.869824
$ws.Range("N4").Value = 2.609472
$ws.Range("O4").Value = 0.008478546007235204
$ws.Range("P4").Value = 0.008478546007235202
$ws.Range("Q4").Value = 0.025640671872
$ws.Range("R4").Value = 0.230766046848
$ws.Range("S4").Value = 0.001301785047890936
$ws.Range("T4").Value = 0.001301785047890935

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Fgf15"
$ws.Range("C5").Value = "Fgfr1"
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.029478
$ws.Range("H5").Value = 0.088434
$ws.Range("I5").Value = 0.1535387136874709
$ws.Range("J5").Value = 0.1535387136874709
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.6143926666666666
$ws.Range("N5").Value = 1.843178
$ws.Range("O5").Value = 0.005988747713148011
$ws.Range("P5").Value = 0.005988747713148009
$ws.Range("Q5").Value = 0.018111067028
$ws.Range("R5").Value = 0.162999603252
$ws.Range("S5").Value = 0.0009195046204755286
$ws.Range("T5").Value = 0.0009195046204755284

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Fgf15"
$ws.Range("C6").Value = "Fgfr1"
$ws.Range("D6").Value = "sCs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.029478
$ws.Range("H6").Value = 0.088434
$ws.Range("I6").Value = 0.1535387136874709
$ws.Range("J6").Value = 0.1535387136874709
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 16.20207433333333
$ws.Range("N6").Value = 48.606223
$ws.Range("O6").Value = 0.1579285380120706
$ws.Range("P6").Value = 0.1579285380120706
$ws.Range("Q6").Value = 0.477604747198
$ws.Range("R6").Value = 4.298442724782
$ws.Range("S6").Value = 0.02424814458091618
$ws.Range("T6").Value = 0.02424814458091617

# Row 7
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Fgf15"
$ws.Range("C7").Value = "Fgfr1"
$ws.Range("D7").Value = "ECs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.1625126666666667
$ws.Range("H7").Value = 0.487538
$ws.Range("I7").Value = 0.846461286312529
$ws.Range("J7").Value = 0.846461286312529
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 10.25548833333333
$ws.Range("N7").Value = 30.766465
$ws.Range("O7").Value = 0.09996462463766297
$ws.Range("P7").Value = 0.09996462463766294
$ws.Range("Q7").Value = 1.666646757018889
$ws.Range("R7").Value = 14.99982081317
$ws.Range("S7").Value = 0.08461618475654532
$ws.Range("T7").Value = 0.08461618475654531

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Fgf15"
$ws.Range("C8").Value = "Fgfr1"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.1625126666666667
$ws.Range("H8").Value = 0.487538
$ws.Range("I8").Value = 0.846461286312529
$ws.Range("J8").Value = 0.846461286312529
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 74.64939600000001
$ws.Range("N8").Value = 223.948188
$ws.Range("O8").Value = 0.7276395436298834
$ws.Range("P8").Value = 0.7276395436298831
$ws.Range("Q8").Value = 12.131472409016
$ws.Range("R8").Value = 109.183251681144
$ws.Range("S8").Value = 0.6159187040728127
$ws.Range("T8").Value = 0.6159187040728125

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Fgf15"
$ws.Range("C9").Value = "Fgfr1"
$ws.Range("D9").Value = "M1"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.1625126666666667
$ws.Range("H9").Value = 0.487538
$ws.Range("I9").Value = 0.846461286312529
$ws.Range("J9").Value = 0.846461286312529
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.869824
$ws.Range("N9").Value = 2.609472
$ws.Range("O9").Value = 0.008478546007235204
$ws.Range("P9").Value = 0.008478546007235202
$ws.Range("Q9").Value = 0.1413574177706667
$ws.Range("R9").Value = 1.272216759936
$ws.Range("S9").Value = 0.007176760959344268
$ws.Range("T9").Value = 0.007176760959344266

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Fgf15"
$ws.Range("C10").Value = "Fgfr1"
$ws.Range("D10").Value = "M2"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.1625126666666667
$ws.Range("H10").Value = 0.487538
$ws.Range("I10").Value = 0.846461286312529
$ws.Range("J10").Value = 0.846461286312529
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.6143926666666666
$ws.Range("N10").Value = 1.843178
$ws.Range("O10").Value = 0.005988747713148011
$ws.Range("P10").Value = 0.005988747713148009
$ws.Range("Q10").Value = 0.09984659064044443
$ws.Range("R10").Value = 0.8986193157639999
$ws.Range("S10").Value = 0.005069243092672482
$ws.Range("T10").Value = 0.005069243092672481

# Row 11
$ws.Range("A11").Value = "sCs"
$ws.Range("B11").Value = "Fgf15"
$ws.Range("C11").Value = "Fgfr1"
$ws.Range("D11").Value = "sCs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.1625126666666667
$ws.Range("H11").Value = 0.487538
$ws.Range("I11").Value = 0.846461286312529
$ws.Range("J11").Value = 0.846461286312529
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 16.20207433333333
$ws.Range("N11").Value = 48.606223
$ws.Range("O11").Value = 0.1579285380120706
$ws.Range("P11").Value = 0.1579285380120706
$ws.Range("Q11").Value = 2.633042305441555
$ws.Range("R11").Value = 23.697380748974
$ws.Range("S11").Value = 0.1336803934311544
$ws.Range("T11").Value = 0.1336803934311544
